$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) First occurrence of "Usará o novo sistema BW_SOC, para:"
#    -> split into 3 runs: "Usará o novo sistema " / "J4Work" / ", para:"
# ------------------------------------------------------------------
$rng1 = $d.Content
[void]$rng1.Find.Execute("BW_SOC")
$rng1.Text = ""
$pos1 = $rng1.Start

$r1 = $d.Range($pos1, $pos1)
$r1.InsertAfter("J4Work")

# Flip the color away and back to force a run boundary without
# changing the visible formatting (keeps rPr identical to neighbors).
$r1b = $d.Range($pos1, $pos1 + 6)
$r1b.Font.Color = 255
$r1b.Font.Color = 0

# ------------------------------------------------------------------
# 2) Second occurrence of "Usará o novo sistema BW_SOC, para:"
#    -> split into 4 runs: "Usará o novo sistema " / "J" / "4Work" / ", para:"
# ------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Start = $pos1 + 6
[void]$rng2.Find.Execute("BW_SOC")
$rng2.Text = ""
$pos2 = $rng2.Start

$r2 = $d.Range($pos2, $pos2)
$r2.InsertAfter("J")
$r3 = $d.Range($pos2 + 1, $pos2 + 1)
$r3.InsertAfter("4Work")

# Apply the color-flip trick right-to-left so earlier splits aren't
# re-coalesced by a later insertion/format operation.
$r3b = $d.Range($pos2 + 1, $pos2 + 6)
$r3b.Font.Color = 255
$r3b.Font.Color = 0

$r2b = $d.Range($pos2, $pos2 + 1)
$r2b.Font.Color = 255
$r2b.Font.Color = 0

# ------------------------------------------------------------------
# 3) "Fazer a gestão de cadastros, compras e vendas."
#    -> "Fazer a gestão de cadastros, compra e venda."
# ------------------------------------------------------------------
[void]$d.Content.Find.Execute("Fazer a gestão de cadastros, compras e vendas.", $true, $false, $false, $false, $false, $true, 1, $false, "Fazer a gestão de cadastros, compra e venda.", 2)

# ------------------------------------------------------------------
# 4) "Organizar o mapeamento de atividades financeiras."
#    -> "Organizar o mapeamento de atividade financeira." split into
#       two runs with the _GoBack bookmark re-inserted right before
#       the final period.
# ------------------------------------------------------------------
$rng4 = $d.Content
[void]$rng4.Find.Execute("atividades financeiras")
$rng4.Text = "atividade financeira"
$bmPos = $rng4.End
$bmRange = $d.Range($bmPos, $bmPos)

# ------------------------------------------------------------------
# 5) Move the _GoBack bookmark: remove it from the final empty
#    paragraph and re-add it right before the period above.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $bmRange)
